# "updated docs for arial 11"
# The whole body of this lesson-plan document is set in Arial; bump every
# run (and every paragraph mark, i.e. the pilcrow's own run properties) from
# the implicit 12pt down to an explicit 11pt (half-points: 22) for both the
# ASCII/Western size (w:sz) and the complex-script size (w:szCs).

$d = $word.ActiveDocument

# Run-level formatting: Font.Size -> w:sz, Font.SizeBi -> w:szCs.
# Touching the whole story (Content) stamps every <w:r><w:rPr> in the body.
$d.Content.Font.Size = 11
$d.Content.Font.SizeBi = 11

# The paragraph mark itself (the run properties living in <w:pPr><w:rPr>)
# is only picked up when the *paragraph's own* Range (which includes the
# end-of-paragraph mark) has its font touched, so walk every paragraph too.
foreach ($p in $d.Paragraphs) {
    $p.Range.Font.Size = 11
    $p.Range.Font.SizeBi = 11
}

# Also flip on hanging/overflow punctuation for the base "Normal" style
# (w:overflowPunct false -> true in the Normal/style0 pPr).
$normal = $d.Styles("Normal")
$normal.ParagraphFormat.HangingPunctuation = $true
